# The document's single paragraph/run originally reads "askdjls dasd".
# This commit just corrects/replaces that text with "ksj dfdfsdfs;d fs";
# every other bit of paragraph/run formatting (font, size, color,
# paragraph spacing/indent/justification, section properties, ...)
# stays exactly as it was.

$d = $word.ActiveDocument

$oldText = "askdjls dasd"
$newText = "ksj dfdfsdfs;d fs"

# Find the range that holds the old text and overwrite just its text
# (instead of a blind Find/Replace) so the run keeps its existing
# formatting/properties untouched.
$rng = $d.Content
$rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true,
                   1, $false, "", 0) | Out-Null
$rng.Text = $newText
